# A new weekly price record is inserted as row 67. Every existing data
# row from 67 through 204 shifts down by one (row 204's old data becomes
# the brand-new row 205), and the sheet dimension grows to A1:R205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the existing rows 67..204 (columns A..R) BEFORE any writes,
#    so the shift-down doesn't clobber data we still need to read.
$srcRange = $ws.Range("A67:R204")
$srcValues = $srcRange.Value2

$firstShiftRow = 68
$srcRowCount = 204 - 67 + 1

# 2) Shift rows 67..204 down into rows 68..205 (row 205 is brand new).
for ($i = 0; $i -lt $srcRowCount; $i++) {
    $destRow = $firstShiftRow + $i
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcValues[$i + 1, $c]
    }
}

# The new row 205 needs the same date-time number format used by the rest
# of column D so it renders/serialises identically to its neighbours.
$ws.Range("D205").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 3) Write the brand-new record into row 67.
$ws.Cells.Item(67, 1).Value2 = 4
$ws.Cells.Item(67, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(67, 3).Value2 = "Los Lagos"
$ws.Cells.Item(67, 4).Value2 = 44662
$ws.Cells.Item(67, 5).Value2 = 10
$ws.Cells.Item(67, 6).Value2 = 100112039
$ws.Cells.Item(67, 7).Value2 = "Ciboulette"
$ws.Cells.Item(67, 8).Value2 = "Sin especificar"
$ws.Cells.Item(67, 9).Value2 = "Primera"
$ws.Cells.Item(67, 10).Value2 = 40
$ws.Cells.Item(67, 11).Value2 = 7000
$ws.Cells.Item(67, 12).Value2 = 7000
$ws.Cells.Item(67, 13).Value2 = 7000
$ws.Cells.Item(67, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(67, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(67, 16).Value2 = 2333
$ws.Cells.Item(67, 17).Value2 = 3
$ws.Cells.Item(67, 18).Value2 = "Hortaliza"
